$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1250
$ws.Range("I70").Value = 1300
$ws.Range("J70").Value = 1200
$ws.Range("K70").Value = 3900
$ws.Range("L70").Value = 3600
$ws.Range("M70").Value = -3630
$ws.Range("N70").Value = -4140
$ws.Range("H73").Value = 1250
$ws.Range("I73").Value = 1300
$ws.Range("J73").Value = 1200
$ws.Range("K73").Value = 3900
$ws.Range("L73").Value = 3600
$ws.Range("M73").Value = -2964
$ws.Range("N73").Value = -5472
$ws.Range("H75").Value = 33000
$ws.Range("J75").Value = 33000
$ws.Range("L75").Value = 33000
$ws.Range("N75").Value = -34872
$ws.Range("H78").Value = 33000
$ws.Range("J78").Value = 33000
$ws.Range("L78").Value = 99000
$ws.Range("N78").Value = -108360
$ws.Range("H100").Value = 2312.7856
$ws.Range("I100").Value = 1942.1111
$ws.Range("J100").Value = 2980
$ws.Range("K100").Value = 1942.1111
$ws.Range("L100").Value = 2980
$ws.Range("M100").Value = -1401.1111
$ws.Range("N100").Value = -4062
$ws.Range("H135").Value = 952.04877
$ws.Range("I135").Value = 499.9375
$ws.Range("K135").Value = 4499.4375
$ws.Range("M135").Value = -1964.4375
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 1690.625
$ws.Range("I137").Value = 1512.5
$ws.Range("J137").Value = 2225
$ws.Range("K137").Value = 4537.5
$ws.Range("L137").Value = 6675
$ws.Range("M137").Value = -1987.5
$ws.Range("N137").Value = -11775
$ws.Range("H138").Value = 2455.506
$ws.Range("I138").Value = 2210.182
$ws.Range("J138").Value = 2546.9832
$ws.Range("K138").Value = 6630.545999999999
$ws.Range("L138").Value = 7640.9496
$ws.Range("M138").Value = -1490.545999999999
$ws.Range("N138").Value = -17920.9496
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 823915.2
$ws.Range("I32").Value = 985784.5
$ws.Range("K32").Value = 985784.5
$ws.Range("M32").Value = -985497.5
$ws.Range("N32").ClearContents()
$ws.Range("H123").Value = 34999.5
$ws.Range("J123").Value = 34999.5
$ws.Range("L123").Value = 34999.5
$ws.Range("N123").Value = -44799.5
$ws.Range("H133").Value = 27232
$ws.Range("J133").Value = 27232
$ws.Range("L133").Value = 27232
$ws.Range("N133").Value = -32292
$ws.Range("H134").Value = 50419
$ws.Range("J134").Value = 50419
$ws.Range("L134").Value = 50419
$ws.Range("N134").Value = -60559
$ws.Range("H135").Value = 24740.584
$ws.Range("J135").Value = 24740.584
$ws.Range("L135").Value = 24740.584
$ws.Range("N135").Value = -34880.584
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 2032.95
$ws.Range("I80").Value = 2387.889
$ws.Range("J80").Value = 1742.5454
$ws.Range("K80").Value = 2387.889
$ws.Range("L80").Value = 1742.5454
$ws.Range("M80").Value = -1389.889
$ws.Range("N80").Value = -3738.5454
$ws.Range("H83").Value = 2032.95
$ws.Range("I83").Value = 2387.889
$ws.Range("J83").Value = 1742.5454
$ws.Range("K83").Value = 11939.445
$ws.Range("L83").Value = 8712.726999999999
$ws.Range("M83").Value = -6947.445
$ws.Range("N83").Value = -18696.727
$ws.Range("H134").Value = 3920.697
$ws.Range("I134").Value = 3318.25
$ws.Range("J134").Value = 4847.5386
$ws.Range("K134").Value = 9954.75
$ws.Range("L134").Value = 14542.6158
$ws.Range("M134").Value = -7419.75
$ws.Range("N134").Value = -19612.6158
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 296.66666
$ws.Range("I7").Value = 316
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 316
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -203
$ws.Range("N7").Value = -426
$ws.Range("H22").Value = 100000420
$ws.Range("I22").Value = 284.2
$ws.Range("J22").Value = 200000540
$ws.Range("K22").Value = 284.2
$ws.Range("L22").Value = 200000540
$ws.Range("M22").Value = 65.80000000000001
$ws.Range("N22").Value = -200001240
$ws.Range("H31").Value = 4401.912
$ws.Range("I31").Value = 1399.091
$ws.Range("J31").Value = 6289.4
$ws.Range("K31").Value = 1399.091
$ws.Range("L31").Value = 6289.4
$ws.Range("M31").Value = -1104.091
$ws.Range("N31").Value = -6879.4
$ws.Range("H34").Value = 4401.912
$ws.Range("I34").Value = 1399.091
$ws.Range("J34").Value = 6289.4
$ws.Range("K34").Value = 1399.091
$ws.Range("L34").Value = 6289.4
$ws.Range("M34").Value = -1197.091
$ws.Range("N34").Value = -6693.4
$ws.Range("H119").Value = 42500
$ws.Range("J119").Value = 42500
$ws.Range("L119").Value = 42500
$ws.Range("N119").Value = -52176
$ws.Range("H130").Value = 58636.363
$ws.Range("J130").Value = 58636.363
$ws.Range("L130").Value = 58636.363
$ws.Range("N130").Value = -68676.363
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 129.375
$ws.Range("I12").Value = 20
$ws.Range("J12").Value = 165.83333
$ws.Range("K12").Value = 60
$ws.Range("L12").Value = 497.49999
$ws.Range("M12").Value = 113
$ws.Range("N12").Value = -843.49999
$ws.Range("H33").Value = 66866.87
$ws.Range("I33").Value = 8541.916999999999
$ws.Range("J33").Value = 300166.66
$ws.Range("K33").Value = 51251.50199999999
$ws.Range("L33").Value = 1800999.96
$ws.Range("M33").Value = -50968.50199999999
$ws.Range("N33").Value = -1801565.96
$ws.Range("H113").Value = 537.5161000000001
$ws.Range("I113").Value = 545.619
$ws.Range("J113").Value = 520.5
$ws.Range("K113").Value = 1636.857
$ws.Range("L113").Value = 1561.5
$ws.Range("M113").Value = 533.143
$ws.Range("N113").Value = -5901.5
$ws.Range("H137").Value = 9268220
$ws.Range("I137").Value = 33354634
$ws.Range("K137").Value = 100063902
$ws.Range("M137").Value = -100058802
$ws.Range("N137").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 1985.8889
$ws.Range("I122").Value = 1410.4286
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4231.2858
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1781.2858
$ws.Range("N122").Value = -16900
$ws.Range("H123").Value = 10309.75
$ws.Range("J123").Value = 10309.75
$ws.Range("L123").Value = 10309.75
$ws.Range("N123").Value = -15209.75
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 40326
$ws.Range("J133").Value = 40326
$ws.Range("L133").Value = 40326
$ws.Range("N133").Value = -45386
$ws.Range("H136").Value = 2689194.8
$ws.Range("I136").Value = 968.0784
$ws.Range("J136").Value = 15152792
$ws.Range("K136").Value = 2904.2352
$ws.Range("L136").Value = 45458376
$ws.Range("M136").Value = -354.2352000000001
$ws.Range("N136").Value = -45463476
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H122").Value = 2272.658
$ws.Range("I122").Value = 2177.853
$ws.Range("J122").Value = 3078.5
$ws.Range("K122").Value = 6533.559
$ws.Range("L122").Value = 9235.5
$ws.Range("M122").Value = -4083.559
$ws.Range("N122").Value = -14135.5
$ws.Range("H123").Value = 27369
$ws.Range("J123").Value = 37685.6
$ws.Range("L123").Value = 37685.6
$ws.Range("N123").Value = -47485.6
